# Automation test case written for OMS-777, OMS-3046, OMS-3470
# Adds a "Vendor Part Number" column to Sheet1 and a new Sheet2 with
# vendor/owner lookup data.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Sheet1: new column K "Vendor Part Number" (same bold/Arial header style
# as the other header cells) ---
$ws1.Range("K1").Font.Bold = $true
$ws1.Range("K1").Font.Name = "Arial"
$ws1.Range("K1").Value = "Vendor Part Number"

# --- Add Sheet2 right after Sheet1 ---
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Sheet2"

$ws2.Range("A1").Value = "OMS-3257"
$ws2.Range("A2").Value = "OMS-3258"
$ws2.Range("B2").Value = "Robin joseph, gajendra M"
$ws2.Range("B1").Value = "Darakshan"
$ws2.Range("B2").Select()

# Leave Sheet1 as the active/selected sheet, with K3 selected
$ws1.Activate()
$ws1.Range("K3").Select()
